# Update odds data in Sheet1 (FlashScore weekly games export).
# Values below reflect the latest refresh of betting odds for the
# matches on rows 4, 6, 8, 9, 10 and 11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 8
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 3

# Row 6
$ws.Range("G6").Value = 1.36
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = 8
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 13
$ws.Range("Q6").Value = 1.67
$ws.Range("R6").Value = 2.15
$ws.Range("U6").Value = 1.95
$ws.Range("V6").Value = 1.8
$ws.Range("Y6").Value = 8.5
$ws.Range("Z6").Value = 9
$ws.Range("AC6").Value = 13
$ws.Range("AD6").Value = 9.5
$ws.Range("AF6").Value = 51
$ws.Range("AI6").Value = 23
$ws.Range("AM6").Value = 351
$ws.Range("AN6").Value = 3.4

# Row 8
$ws.Range("K8").Value = 2.05
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 7.5
$ws.Range("O8").Value = 1.4
$ws.Range("P8").Value = 2.75
$ws.Range("AC8").Value = 7.5
$ws.Range("AP8").Value = 23
$ws.Range("AS8").Value = 201

# Row 9
$ws.Range("G9").Value = 2.38
$ws.Range("I9").Value = 3.2
$ws.Range("J9").Value = 3.2
$ws.Range("L9").Value = 4
$ws.Range("M9").Value = 1.11
$ws.Range("N9").Value = 6.5
$ws.Range("S9").Value = 1.53
$ws.Range("T9").Value = 2.38
$ws.Range("W9").Value = 6.5
$ws.Range("X9").Value = 10
$ws.Range("Z9").Value = 23
$ws.Range("AA9").Value = 23
$ws.Range("AE9").Value = 17
$ws.Range("AG9").Value = 7.5
$ws.Range("AI9").Value = 12
$ws.Range("AJ9").Value = 34
$ws.Range("AK9").Value = 29
$ws.Range("AN9").Value = 4.33
$ws.Range("AO9").Value = 15
$ws.Range("AT9").Value = 2.38
$ws.Range("AX9").Value = 19

# Row 10
$ws.Range("G10").Value = 1.57
$ws.Range("H10").Value = 4.33
$ws.Range("I10").Value = 5
$ws.Range("J10").Value = 2.1
$ws.Range("K10").Value = 2.5
$ws.Range("L10").Value = 5
$ws.Range("N10").Value = 17
$ws.Range("X10").Value = 8.5
$ws.Range("Z10").Value = 12
$ws.Range("AC10").Value = 17
$ws.Range("AD10").Value = 8.5
$ws.Range("AG10").Value = 17
$ws.Range("AH10").Value = 29
$ws.Range("AO10").Value = 8
$ws.Range("AW10").Value = 7
$ws.Range("AX10").Value = 26
$ws.Range("AY10").Value = 29

# Row 11
$ws.Range("G11").Value = 5.25
$ws.Range("H11").Value = 4
$ws.Range("I11").Value = 1.6
$ws.Range("K11").Value = 2.3
$ws.Range("L11").Value = 2.2
$ws.Range("Z11").Value = 51
$ws.Range("AC11").Value = 13
$ws.Range("AD11").Value = 7.5
$ws.Range("AE11").Value = 15
$ws.Range("AH11").Value = 8
$ws.Range("AJ11").Value = 12
$ws.Range("AK11").Value = 13
$ws.Range("AM11").Value = 201
$ws.Range("AN11").Value = 7
$ws.Range("AO11").Value = 26
$ws.Range("AP11").Value = 29
$ws.Range("AQ11").Value = 81
$ws.Range("AR11").Value = 101
$ws.Range("AU11").Value = 8
$ws.Range("AX11").Value = 8
$ws.Range("AY11").Value = 17
